$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide "CONTROL BLOCK" (currently slide 1), Subtitle shape: merge the
#    "seus " + "proble" + "mas" runs into a single "seus problemas" run.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(3)
$tr1 = $subtitle.TextFrame.TextRange
$full1 = $tr1.Text
$idx1 = $full1.IndexOf("seus ")
$target1 = "seus problemas"
$sub1 = $tr1.Characters($idx1 + 1, $target1.Length)
$sub1.Text = $target1

# ---------------------------------------------------------------------------
# 2) "Control bus" overview slide (currently slide 2), content placeholder:
#    merge the " " + "e " runs (between "ônibus" and "sua") into " e ".
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$content2 = $s2.Shapes.Item(2)
$tr2 = $content2.TextFrame.TextRange
$full2 = $tr2.Text
$idx2 = $full2.IndexOf("nibus e sua")
$target2 = " e "
$sub2 = $tr2.Characters($idx2 + 1 + 5, $target2.Length)
$sub2.Text = $target2

# ---------------------------------------------------------------------------
# 3) "Quanto custa" slide (currently slide 4), first content placeholder:
#    reposition the shape and justify its paragraph.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$content4 = $s4.Shapes.Item(2)
$content4.Left = 66.97142028808594
$content4.Top = 161.88575744628906
$content4.TextFrame.TextRange.ParagraphFormat.Alignment = 4

# ---------------------------------------------------------------------------
# 4) Move the "Control bus" overview slide from position 2 to position 5
#    (it now follows the "Quanto custa" / empty-placeholder slides).
# ---------------------------------------------------------------------------
$s2.MoveTo(5)
